$wb = $excel.ActiveWorkbook

# Sheets (by index, per workbook.xml order):
#   1 = SAD-Code
#   2 = SAD-SAM
#   3 = SAM-Code
#   4 = Projects
$wsSadCode = $wb.Worksheets.Item(1)
$wsSadSam  = $wb.Worksheets.Item(2)
$wsSamCode = $wb.Worksheets.Item(3)
$wsProjects = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# Update gold standard confusion-matrix counts for BigBlueButton (row 4)
# on the Projects sheet. Column C = SAD-Code, D = SAD-SAM, E = SAM-Code.
# ---------------------------------------------------------------------
$wsProjects.Range("C4").Value = 1473
$wsProjects.Range("E4").Value = 730

# ---------------------------------------------------------------------
# Update results for BigBlueButton (row 4) on the SAD-Code sheet.
# ---------------------------------------------------------------------
$wsSadCode.Range("C4").Value = 0.77
$wsSadCode.Range("D4").Value = 0.91
$wsSadCode.Range("E4").Value = 0.84

# ---------------------------------------------------------------------
# Update results for BigBlueButton (row 4) on the SAM-Code sheet.
# ---------------------------------------------------------------------
$wsSamCode.Range("C4").Value = 0.94
$wsSamCode.Range("D4").Value = 0.96
$wsSamCode.Range("E4").Value = 0.95
$wsSamCode.Range("G4").Value = 1
$wsSamCode.Range("H4").Value = 0.95

# ---------------------------------------------------------------------
# Update the selected / active cell on each sheet to match the saved
# workbook view state.
# ---------------------------------------------------------------------
$wsSadCode.Activate() | Out-Null
$wsSadCode.Range("E4").Select() | Out-Null

$wsSadSam.Activate() | Out-Null
$wsSadSam.Range("C4").Select() | Out-Null

$wsSamCode.Activate() | Out-Null
$wsSamCode.Range("C4").Select() | Out-Null

$wsProjects.Activate() | Out-Null
$wsProjects.Range("C5").Select() | Out-Null

# Re-activate the sheet that was active originally.
$wsSadCode.Activate() | Out-Null
